$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 8.723478434330044
$ws.Range("I10").Value = 2.080435302990134
$ws.Range("K10").Value = 0.0
$ws.Range("P10").Value = 4.572825479844843
$ws.Range("R10").Value = 2.409837559296905
$ws.Range("S10").Value = 2.080435302990134
$ws.Range("T10").Value = 2.3196853628339995
$ws.Range("U10").Value = 1.8789429385525764
$ws.Range("V10").Value = 0.0
$ws.Range("W10").Value = 2.1671201072813897
$ws.Range("X10").Value = 4.137128678467276
$ws.Range("Y10").Value = 6.142611122979786
$ws.Range("Z10").Value = 2.4271745201551562
$ws.Range("AA10").Value = 0.0
$ws.Range("AI10").Value = 6.522921419289059
$ws.Range("AL10").Value = 2.888337678984636
$ws.Range("AN10").Value = 2.1948592446545914
$ws.Range("AO10").Value = 2.4271745201551562
$ws.Range("AR10").Value = 4.586203633873944
$ws.Range("AS10").Value = 2.8665872266164922
$ws.Range("AT10").Value = 1.3414511924000259
$ws.Range("AU10").Value = 2.080435302990134
$ws.Range("AV10").Value = 0.0
$ws.Range("AY10").Value = 4.503829698515615
$ws.Range("BA10").Value = 8.850624004490722
